$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.206.62"
$ws.Range("E2").Value = "  +3.02%  "

$ws.Range("D3").Value = "2.058.98"
$ws.Range("E3").Value = "  +2.37%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.09"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.87"
$ws.Range("E7").Value = "  +8.48%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +3.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("E10").Value = "  +4.62%  "

$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.363.92"
$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.73"
$ws.Range("E13").Value = "  +5.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.22"
$ws.Range("E14").Value = "  +7.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").Value = "2.078.20"
$ws.Range("E17").Value = "  +2.31%  "

$ws.Range("D18").Value = "38.063.06"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("E21").Value = "  +3.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.49"
$ws.Range("E22").Value = "  +0.77%  "

$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("E25").Value = "  +4.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.54"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("E27").Value = "  +4.13%  "

$ws.Range("E28").Value = "  +7.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.06"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("E30").Value = "  +2.26%  "

$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("E32").Value = "  +3.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +3.21%  "

$ws.Range("E34").Value = "  +10.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0606"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.32"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.18"
$ws.Range("E37").Value = "  +15.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +5.93%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").Value = "1.525.30"
$ws.Range("E40").Value = "  +4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.12"
$ws.Range("E41").Value = "  +7.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.16"
$ws.Range("E42").Value = "  +3.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0217"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.86"
$ws.Range("E44").Value = "  +3.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.06"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.98"
$ws.Range("E49").Value = "  +3.17%  "

$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("D51").Value = "2.251.71"
$ws.Range("E51").Value = "  +2.42%  "
